# Logboek voetbal.xlsx update
# - resize book window view
# - rewrite several "E" column remarks (strip parentheses / "door")
# - add a new block of logboek entries for May (rows 26-38)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E remarks that were reworded (parentheses / "door" removed) ---
$ws.Range("E7").Value  = "ontwerpschets iedereen, samenwerkingscontract en bereikbaarheidslijst Corwin"
$ws.Range("E10").Value = "Conventierapport door Sven"
$ws.Range("E11").Value = "PvA Yannick, Rolverdeling + GitHub accounts iedereen"
$ws.Range("E12").Value = "Voorkant Steven"
$ws.Range("E18").Value = "Logboek Steven en Yannick, verbeteren PvA Corwin "

# --- Row 26 used to be the "6-11 apr" placeholder row; it now becomes the
#     "6-11 mei" header row that kicks off the new May entries. ---
$ws.Range("C26").Value = "6-11 mei"
$ws.Range("D26").Value = "meivakantie, prototype gemaakt"
$ws.Range("E26").Value = "Prototype Yannick"

# --- Row 28: date moves from 12 apr to 12 mei, text updated ---
$ws.Range("C28").NumberFormat = "d-mmm"
$ws.Range("C28").Value = 41771
$ws.Range("D28").Value = "Prototype verder gemaakt, rapport database > programma, database onderzoek verslag, flowchart"
$ws.Range("E28").Value = "Prototype Yannick en Corwin, rapport database naar programma Sven, onderzoek database verslag Steven"

# --- Row 30: date moves from 13 apr to 13 mei, new text added ---
$ws.Range("C30").NumberFormat = "d-mmm"
$ws.Range("C30").Value = 41772
$ws.Range("D30").Value = "Onderzoek database verslag, "
$ws.Range("E30").Value = "verslag Steven, "

# --- New row 32 (19 mei) ---
$ws.Range("C32").NumberFormat = "d-mmm"
$ws.Range("C32").Value = 41778
$ws.Range("D32").Value = "Onderzoek database verslag verder, PvA aangepast, GitHub volledig werkend"
$ws.Range("E32").Value = "database verslag Steven, PvA Yannick, GitHub iedereen"

# --- New rows 34, 36, 38: just dates (20, 26, 27 mei) ---
$ws.Range("C34").NumberFormat = "d-mmm"
$ws.Range("C34").Value = 41779

$ws.Range("C36").NumberFormat = "d-mmm"
$ws.Range("C36").Value = 41785

$ws.Range("C38").NumberFormat = "d-mmm"
$ws.Range("C38").Value = 41786

# --- column widths (bestFit-ish; engine quantizes to its own grid so these
#     are the closest achievable inputs to the target stored widths) ---
$ws.Columns.Item(3).ColumnWidth = 7.6
$ws.Columns.Item(4).ColumnWidth = 89.0
$ws.Columns.Item(5).ColumnWidth = 96.3

# --- selection moves to the next empty row after the new data ---
$ws.Range("C39").Select() | Out-Null

# --- resize the workbook window (maximized-ish view) ---
$excel.ActiveWindow.WindowState = -4137
